$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stable style "donors" that are not themselves modified by this edit:
#   H8  -> fillId=3 ("highlight" fill, style index 5)
#   F20 -> fillId=2 (the other fill, style index 2)

# Row 11: add F11, new value, highlighted fill (s=5)
$ws.Range("F11").Value = "Prep load"
$ws.Range("H8").Copy()
$ws.Range("F11").PasteSpecial(-4122)

# Row 12: F12 keeps its text but switches fill to highlighted (s=5); add H12 highlighted (s=5)
$ws.Range("H8").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("H12").Value = "Rearrange"
$ws.Range("H8").Copy()
$ws.Range("H12").PasteSpecial(-4122)

# Row 13: F13 keeps its text but switches fill to highlighted (s=5)
$ws.Range("H8").Copy()
$ws.Range("F13").PasteSpecial(-4122)

# Row 14: add H14 with default (no fill) style
$ws.Range("H14").Value = "Prep chambers 11:30"

# Row 15: H15 keeps its text but switches fill to the other fill (s=2)
$ws.Range("F20").Copy()
$ws.Range("H15").PasteSpecial(-4122)

# Row 16: add G16 with default style; H16 keeps its text but switches fill to the other fill (s=2)
$ws.Range("G16").Value = "Prep load "
$ws.Range("F20").Copy()
$ws.Range("H16").PasteSpecial(-4122)

# Row 17: add F17 with default style; H17 (empty) switches fill to the other fill (s=2)
$ws.Range("F17").Value = "Rearrange"
$ws.Range("F20").Copy()
$ws.Range("H17").PasteSpecial(-4122)

# Row 19: add F19 with default style
$ws.Range("F19").Value = "Prep chambers 14:00"

# Row 21: add G21 with default style
$ws.Range("G21").Value = "Rearrange"

# Row 24: add G24 with default style
$ws.Range("G24").Value = "Prep chambers 16:30"

# --- Sheet view changes: scroll so column C is the left-most visible column, and select H14 ---
$ws.Range("H14").Select()
$excel.ActiveWindow.ScrollColumn = 3
